# Dodajanje diplomaticnega prepisa Kosijeve in Sefove pesmarice
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 (Tjasa Miholic / NUK MS 1485-M13 / f. 97-144): status finished ->
# mark "koncano" and highlight it green to flag the newly finished transcript.
$ws.Range("D4").Value = "končano"
$ws.Range("D4").Interior.Color = 5296274   # RGB(146,208,80) == FF92D050

# Row 5 (Akos Doncec / Martjanska II): no longer has a "v delu" status value.
$ws.Range("D5").Value = ""

# New rows for the diplomatic transcription of the Kosi songbook
# (Markiševska manuscript), written in an order that keeps the new shared
# strings appended the same way the author's workbook did.
$ws.Range("C11").Value = "f. 1-121"
$ws.Range("B11").Value = "Markiševska"
$ws.Range("A11").Value = "Špela Kovačič"
$ws.Range("D11").Value = "končano"
$ws.Range("E11").Value = "da"
$ws.Range("D2").Copy()
$ws.Range("D11").PasteSpecial(-4122)   # xlPasteFormats - reuse the "koncano" highlight style

$ws.Range("C12").Value = "122-244"
$ws.Range("B12").Value = "Markiševska"
$ws.Range("A12").Value = "Tjaša Miholič"

# Row 7 used to hold a stray "Klementina Kosi" entry with no manuscript data;
# replace it with the Kosijeva pesmarica entry for Maja Lampret.
$ws.Range("B7").Value = "Kosijeva pesmarica"
$ws.Range("A7").Value = "Maja Lampret"

# New row for the diplomatic transcription of the Sef songbook
$ws.Range("A13").Value = "Špeka Kovačič"
$ws.Range("B13").Value = "Šefova"

$ws.Range("C13").Select()
